$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new blank rows before row 74, shifting existing rows 74-180 down to 76-182
$ws.Rows.Item(74).Insert()
$ws.Rows.Item(74).Insert()

# Common (static) field values reused across rows on this sheet
$mercadoId = 7
$mercado = "Terminal Hortofrutícola Agro Chillán"
$region = "Ñuble"
$codreg = 16
$tipo = "Fruta"
$productoId = 100108
$producto = "Tropicales y subtropicales"
$subproductoId = 100108005
$subproducto = "Piña"
$origen = "Ecuador"

# New row 74: Primera, 12-unit box
$ws.Range("A74").Value = $mercadoId
$ws.Range("B74").Value = $mercado
$ws.Range("C74").Value = $region
$ws.Range("D74").Value = 44540
$ws.Range("E74").Value = $codreg
$ws.Range("F74").Value = $tipo
$ws.Range("G74").Value = $productoId
$ws.Range("H74").Value = $producto
$ws.Range("I74").Value = $subproductoId
$ws.Range("J74").Value = $subproducto
$ws.Range("K74").Value = "Caramelo"
$ws.Range("L74").Value = "Primera"
$ws.Range("M74").Value = 60
$ws.Range("N74").Value = 20000
$ws.Range("O74").Value = 21000
$ws.Range("P74").Value = 20500
$ws.Range("Q74").Value = "$/caja 12 unidades"
$ws.Range("R74").Value = $origen
$ws.Range("S74").Value = 1708
$ws.Range("T74").Value = 12

# New row 75: Segunda, 14-unit box
$ws.Range("A75").Value = $mercadoId
$ws.Range("B75").Value = $mercado
$ws.Range("C75").Value = $region
$ws.Range("D75").Value = 44540
$ws.Range("E75").Value = $codreg
$ws.Range("F75").Value = $tipo
$ws.Range("G75").Value = $productoId
$ws.Range("H75").Value = $producto
$ws.Range("I75").Value = $subproductoId
$ws.Range("J75").Value = $subproducto
$ws.Range("K75").Value = "Caramelo"
$ws.Range("L75").Value = "Segunda"
$ws.Range("M75").Value = 60
$ws.Range("N75").Value = 18000
$ws.Range("O75").Value = 19000
$ws.Range("P75").Value = 18500
$ws.Range("Q75").Value = "$/caja 14 unidades"
$ws.Range("R75").Value = $origen
$ws.Range("S75").Value = 1321
$ws.Range("T75").Value = 14
